$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.716.46"
$ws.Range("E2").Value = "  +5.85%  "
$ws.Range("D3").Value = "3.634.38"
$ws.Range("E3").Value = "  +5.83%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "593.06"
$ws.Range("E5").Value = "  +2.43%  "
$ws.Range("D6").Value = "195.16"
$ws.Range("E6").Value = "  +3.93%  "
$ws.Range("D7").Value = "0.645"
$ws.Range("E7").Value = "  +2.83%  "
$ws.Range("D8").Value = "3.627.87"
$ws.Range("E8").Value = "  +5.85%  "
$ws.Range("D10").Value = "0.184"
$ws.Range("E10").Value = "  +8.19%  "
$ws.Range("D11").Value = "0.678"
$ws.Range("E11").Value = "  +5.67%  "
$ws.Range("D12").Value = "58.02"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "0.0000311"
$ws.Range("E13").Value = "  +13.14%  "
$ws.Range("D14").Value = "9.94"
$ws.Range("E14").Value = "  +5.65%  "
$ws.Range("D15").Value = "4.218.60"
$ws.Range("E15").Value = "  +5.99%  "
$ws.Range("D16").Value = "20.46"
$ws.Range("E16").Value = "  +8.37%  "
$ws.Range("D17").Value = "3.632.29"
$ws.Range("E17").Value = "  +5.40%  "
$ws.Range("D18").Value = "70.744.55"
$ws.Range("E18").Value = "  +5.95%  "
$ws.Range("D19").Value = "12.74"
$ws.Range("E19").Value = "  +6.18%  "
$ws.Range("E20").Value = "  +2.87%  "
$ws.Range("D21").Value = "1.06"
$ws.Range("E21").Value = "  +4.47%  "
$ws.Range("D22").Value = "488.48"
$ws.Range("E22").Value = "  +2.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "19.50"
$ws.Range("E23").Value = "  +15.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.20"
$ws.Range("E24").Value = "  -4.19%  "
$ws.Range("E25").Value = "  +3.59%  "
$ws.Range("D26").Value = "91.32"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("E27").Value = "  +7.25%  "
$ws.Range("D28").Value = "11.42"
$ws.Range("E28").Value = "  +4.70%  "
$ws.Range("E29").Value = "  +7.04%  "
$ws.Range("D30").Value = "7.98"
$ws.Range("E30").Value = "  +9.25%  "
$ws.Range("D31").Value = "32.84"
$ws.Range("E31").Value = "  +5.90%  "
$ws.Range("E32").Value = "  +10.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.30"
$ws.Range("E33").Value = "  +5.07%  "
$ws.Range("D34").Value = "66.23"
$ws.Range("E34").Value = "  +2.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "611.50"
$ws.Range("E35").Value = "  +2.05%  "
$ws.Range("D36").Value = "40.52"
$ws.Range("E36").Value = "  +10.10%  "
$ws.Range("D37").Value = "0.0₃0842"
$ws.Range("E37").Value = "  +13.27%  "
$ws.Range("D38").Value = "0.412"
$ws.Range("E38").Value = "  +6.26%  "
$ws.Range("E39").Value = "  +1.63%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").Value = "3.317.17"
$ws.Range("E42").Value = "  +4.10%  "
$ws.Range("D43").Value = "3.18"
$ws.Range("E43").Value = "  +10.26%  "
$ws.Range("B44").Value = "Fetch.AI"
$ws.Range("C44").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D44").Value = "2.88"
$ws.Range("E44").Value = "  +11.75%  "
$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D45").Value = "3.14"
$ws.Range("E45").Value = "  +18.09%  "
$ws.Range("D46").Value = "0.0458"
$ws.Range("E46").Value = "  +6.95%  "
$ws.Range("D47").Value = "9.67"
$ws.Range("E47").Value = "  +12.67%  "
$ws.Range("D48").Value = "3.35"
$ws.Range("E48").Value = "  +3.05%  "
$ws.Range("E49").Value = "  +3.56%  "
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("D51").Value = "0.999"
$ws.Range("E51").Value = "  -0.08%  "
